# Auto-generated edit script: updates computed market-price snapshot values
# across the 8 crafting-class Leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Values are static (no formulas in the workbook) so each cell is written directly.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 33
$ws.Range("H33").Value = 293.33334
$ws.Range("I33").Value = 290
$ws.Range("K33").Value = 290
$ws.Range("M33").Value = -61

# Row 137
$ws.Range("H137").Value = 2327536.2
$ws.Range("I137").Value = 3126411.2
$ws.Range("J137").Value = 3536.6365
$ws.Range("K137").Value = 9379233.600000001
$ws.Range("L137").Value = 10609.9095
$ws.Range("M137").Value = -9376683.600000001
$ws.Range("N137").Value = -15709.9095

$ws = $wb.Worksheets.Item("ARM")

# Row 2
$ws.Range("H2").Value = 1786.2858
$ws.Range("I2").Value = 1752.75
$ws.Range("J2").Value = 1831
$ws.Range("K2").Value = 1752.75
$ws.Range("L2").Value = 1831
$ws.Range("M2").Value = -1639.75
$ws.Range("N2").Value = -2057

# Row 32
$ws.Range("H32").Value = 4032956.5
$ws.Range("I32").Value = 4847652.5
$ws.Range("K32").Value = 4847652.5
$ws.Range("M32").Value = -4847365.5

# Row 116
$ws.Range("H116").Value = 1786.2858
$ws.Range("I116").Value = 1752.75
$ws.Range("J116").Value = 1831
$ws.Range("K116").Value = 1752.75
$ws.Range("L116").Value = 1831
$ws.Range("M116").Value = 541.25
$ws.Range("N116").Value = -6419

# Row 132
$ws.Range("H132").Value = 79140.234
$ws.Range("I132").Value = 57208.277
$ws.Range("J132").Value = 128487.125
$ws.Range("K132").Value = 171624.831
$ws.Range("L132").Value = 385461.375
$ws.Range("M132").Value = -169094.831
$ws.Range("N132").Value = -390521.375

$ws = $wb.Worksheets.Item("BSM")

# Row 3
$ws.Range("H3").Value = 1786.2858
$ws.Range("I3").Value = 1752.75
$ws.Range("J3").Value = 1831
$ws.Range("K3").Value = 1752.75
$ws.Range("L3").Value = 1831
$ws.Range("M3").Value = -1638.75
$ws.Range("N3").Value = -2059

# Row 21
$ws.Range("H21").Value = 15900
$ws.Range("J21").Value = 15900
$ws.Range("L21").Value = 15900
$ws.Range("N21").Value = -16372

$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 9714.617
$ws.Range("I31").Value = 21522.408
$ws.Range("J31").Value = 1938.7561
$ws.Range("K31").Value = 21522.408
$ws.Range("L31").Value = 1938.7561
$ws.Range("M31").Value = -21227.408
$ws.Range("N31").Value = -2528.7561

# Row 34
$ws.Range("H34").Value = 9714.617
$ws.Range("I34").Value = 21522.408
$ws.Range("J34").Value = 1938.7561
$ws.Range("K34").Value = 21522.408
$ws.Range("L34").Value = 1938.7561
$ws.Range("M34").Value = -21320.408
$ws.Range("N34").Value = -2342.7561

# Row 58
$ws.Range("H58").Value = 17896922
$ws.Range("I58").Value = 24288118
$ws.Range("K58").Value = 24288118
$ws.Range("M58").Value = -24287915

# Row 105
$ws.Range("H105").Value = 1303
$ws.Range("I105").Value = 1336.0555
$ws.Range("J105").Value = 1005.5
$ws.Range("K105").Value = 1336.0555
$ws.Range("L105").Value = 1005.5
$ws.Range("M105").Value = 410.9445000000001
$ws.Range("N105").Value = -4499.5

# Row 136
$ws.Range("H136").Value = 17896922
$ws.Range("I136").Value = 24288118
$ws.Range("K136").Value = 72864354
$ws.Range("M136").Value = -72861804

$ws = $wb.Worksheets.Item("CUL")

# Row 5
$ws.Range("H5").Value = 32838.227
$ws.Range("I5").Value = 67140.734
$ws.Range("J5").Value = 679.625
$ws.Range("K5").Value = 201422.202
$ws.Range("L5").Value = 2038.875
$ws.Range("M5").Value = -201310.202
$ws.Range("N5").Value = -2262.875

# Row 80
$ws.Range("H80").Value = 3805.9092
$ws.Range("I80").Value = 2002
$ws.Range("J80").Value = 3986.3
$ws.Range("K80").Value = 6006
$ws.Range("L80").Value = 11958.9
$ws.Range("M80").Value = -5070
$ws.Range("N80").Value = -13830.9

# Row 83
$ws.Range("H83").Value = 3805.9092
$ws.Range("I83").Value = 2002
$ws.Range("J83").Value = 3986.3
$ws.Range("K83").Value = 18018
$ws.Range("L83").Value = 35876.7
$ws.Range("M83").Value = -13338
$ws.Range("N83").Value = -45236.7

# Row 135
$ws.Range("H135").Value = 32838.227
$ws.Range("I135").Value = 67140.734
$ws.Range("J135").Value = 679.625
$ws.Range("K135").Value = 604266.6059999999
$ws.Range("L135").Value = 6116.625
$ws.Range("M135").Value = -601731.6059999999
$ws.Range("N135").Value = -11186.625

$ws = $wb.Worksheets.Item("GSM")

# Row 32
$ws.Range("H32").Value = 28000
$ws.Range("J32").Value = 28000
$ws.Range("L32").Value = 28000
$ws.Range("N32").Value = -28592

# Row 99
$ws.Range("H99").Value = 8823.588
$ws.Range("I99").Value = 4338
$ws.Range("J99").Value = 29756.334
$ws.Range("K99").Value = 4338
$ws.Range("L99").Value = 29756.334
$ws.Range("M99").Value = -2092
$ws.Range("N99").Value = -34248.334

# Row 107
$ws.Range("H107").Value = 1391
$ws.Range("I107").Value = 1054.1428
$ws.Range("J107").Value = 1727.8572
$ws.Range("K107").Value = 1054.1428
$ws.Range("L107").Value = 1727.8572
$ws.Range("M107").Value = 865.8571999999999
$ws.Range("N107").Value = -5567.8572

$ws = $wb.Worksheets.Item("LTW")

# Row 7
$ws.Range("H7").Value = 2563.8572
$ws.Range("I7").Value = 2701.3333
$ws.Range("K7").Value = 2701.3333
$ws.Range("M7").Value = -2589.3333

# Row 40
$ws.Range("H40").Value = 5515.0527
$ws.Range("I40").Value = 5657.2856
$ws.Range("K40").Value = 5657.2856
$ws.Range("M40").Value = -5521.2856

# Row 95
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

# Row 126
$ws.Range("H126").Value = 2563.8572
$ws.Range("I126").Value = 2701.3333
$ws.Range("K126").Value = 8103.999899999999
$ws.Range("M126").Value = -5633.999899999999

$ws = $wb.Worksheets.Item("WVR")

# Row 56
$ws.Range("H56").Value = 15266.667
$ws.Range("I56").Value = 4950
$ws.Range("J56").Value = 35900
$ws.Range("K56").Value = 4950
$ws.Range("L56").Value = 35900
$ws.Range("M56").Value = -4236
$ws.Range("N56").Value = -37328

# Row 57
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()

# Row 97
$ws.Range("H97").Value = 24143
$ws.Range("I97").Value = 10000
$ws.Range("J97").Value = 28857.334
$ws.Range("K97").Value = 10000
$ws.Range("L97").Value = 28857.334
$ws.Range("M97").Value = -9009
$ws.Range("N97").Value = -30839.334

# Row 129
$ws.Range("H129").Value = 28000
$ws.Range("J129").Value = 28000
$ws.Range("L129").Value = 28000
$ws.Range("N129").Value = -38000
